$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2159311196958242
$ws.Cells.Item(2, 4).Value = 0.246097157014546
$ws.Cells.Item(2, 5).Value = 0.1773659978598765
$ws.Cells.Item(2, 6).Value = 1.012249720703693
$ws.Cells.Item(2, 7).Value = 0.5272962014676921
$ws.Cells.Item(2, 8).Value = 0.5812333897015094
$ws.Cells.Item(2, 9).Value = 0.3757315010515043
$ws.Cells.Item(2, 10).Value = 0.1772330055721199
$ws.Cells.Item(2, 14).Value = 1.538965158849294
$ws.Cells.Item(2, 15).Value = 2.183072434228677

$ws.Cells.Item(3, 2).Value = 0.1913949238933697
$ws.Cells.Item(3, 4).Value = 0.2444676174134202
$ws.Cells.Item(3, 5).Value = 0.1742912681868667
$ws.Cells.Item(3, 6).Value = 0.9831260700589439
$ws.Cells.Item(3, 7).Value = 0.4999711970626919
$ws.Cells.Item(3, 8).Value = 0.5723711798586777
$ws.Cells.Item(3, 9).Value = 0.3656238581545637
$ws.Cells.Item(3, 10).Value = 0.1721769985554289
$ws.Cells.Item(3, 14).Value = 1.440536212657179
$ws.Cells.Item(3, 15).Value = 2.106238768013412

$ws.Cells.Item(4, 2).Value = 0.1762869034279646
$ws.Cells.Item(4, 4).Value = 0.2435737801817908
$ws.Cells.Item(4, 5).Value = 0.1724976178116897
$ws.Cells.Item(4, 6).Value = 0.9658221175057946
$ws.Cells.Item(4, 7).Value = 0.4835087855673947
$ws.Cells.Item(4, 8).Value = 0.5672287518405454
$ws.Cells.Item(4, 9).Value = 0.3596127975951688
$ws.Cells.Item(4, 10).Value = 0.1691764068403643
$ws.Cells.Item(4, 14).Value = 1.380350891855613
$ws.Cells.Item(4, 15).Value = 2.060343181411326

$ws.Cells.Item(5, 2).Value = 0.1701201473106835
$ws.Cells.Item(5, 4).Value = 0.2432364075533258
$ws.Cells.Item(5, 5).Value = 0.1717903946035122
$ws.Cells.Item(5, 6).Value = 0.9589159623540269
$ws.Cells.Item(5, 7).Value = 0.4768793569706844
$ws.Cells.Item(5, 8).Value = 0.5652084079914204
$ws.Cells.Item(5, 9).Value = 0.3572128864359954
$ws.Cells.Item(5, 10).Value = 0.167979719279117
$ws.Cells.Item(5, 14).Value = 1.355891405433653
$ws.Cells.Item(5, 15).Value = 2.041962182227564

$ws.Cells.Item(6, 2).Value = 0.1690955695415823
$ws.Cells.Item(6, 4).Value = 0.24318201137892
$ws.Cells.Item(6, 5).Value = 0.171674392540087
$ws.Cells.Item(6, 6).Value = 0.957777978631924
$ws.Cells.Item(6, 7).Value = 0.4757833212346156
$ws.Cells.Item(6, 8).Value = 0.5648774770931766
$ws.Cells.Item(6, 9).Value = 0.3568173988915788
$ws.Cells.Item(6, 10).Value = 0.1677825844652077
$ws.Cells.Item(6, 14).Value = 1.351834048475041
$ws.Cells.Item(6, 15).Value = 2.038929463360034

$ws.Cells.Item(7, 2).Value = 0.1762037765295901
$ws.Cells.Item(7, 4).Value = 0.2435691213927882
$ws.Cells.Item(7, 5).Value = 0.1724879839689208
$ws.Cells.Item(7, 6).Value = 0.9657283902091791
$ws.Cells.Item(7, 7).Value = 0.4834190585289946
$ws.Cells.Item(7, 8).Value = 0.5672012000805466
$ws.Cells.Item(7, 9).Value = 0.3595802297560198
$ws.Cells.Item(7, 10).Value = 0.1691601623408303
$ws.Cells.Item(7, 14).Value = 1.38002074861987
$ws.Cells.Item(7, 15).Value = 2.060093986032228

$ws.Cells.Item(8, 2).Value = 0.2074802658001431
$ws.Cells.Item(8, 4).Value = 0.2455131683178706
$ws.Cells.Item(8, 5).Value = 0.1762862718187961
$ws.Cells.Item(8, 6).Value = 1.002087880873745
$ws.Cells.Item(8, 7).Value = 0.517809018632434
$ws.Cells.Item(8, 8).Value = 0.5781156669493015
$ws.Cells.Item(8, 9).Value = 0.3722062484302811
$ws.Cells.Item(8, 10).Value = 0.1754681284190411
$ws.Cells.Item(8, 14).Value = 1.504976964632107
$ws.Cells.Item(8, 15).Value = 2.156314253543428

$ws.Cells.Item(9, 2).Value = 0.2684501888315936
$ws.Cells.Item(9, 4).Value = 0.2501707657315677
$ws.Cells.Item(9, 5).Value = 0.1844829058625663
$ws.Cells.Item(9, 6).Value = 1.077981530534032
$ws.Cells.Item(9, 7).Value = 0.5877591867547096
$ws.Cells.Item(9, 8).Value = 0.6018910408742215
$ws.Cells.Item(9, 9).Value = 0.3984893586955494
$ws.Cells.Item(9, 10).Value = 0.1886638921441488
$ws.Cells.Item(9, 14).Value = 1.751860150726543
$ws.Cells.Item(9, 15).Value = 2.355184202377018

$ws.Cells.Item(10, 2).Value = 0.3129948669844396
$ws.Cells.Item(10, 4).Value = 0.2541069149559831
$ws.Cells.Item(10, 5).Value = 0.1909625588303214
$ws.Cells.Item(10, 6).Value = 1.136556074637113
$ws.Cells.Item(10, 7).Value = 0.6407026906574345
$ws.Cells.Item(10, 8).Value = 0.6208072618303788
$ws.Cells.Item(10, 9).Value = 0.4186982167615625
$ws.Cells.Item(10, 10).Value = 0.1988664766722508
$ws.Cells.Item(10, 14).Value = 1.93419659562025
$ws.Cells.Item(10, 15).Value = 2.507549276033956

$ws.Cells.Item(11, 2).Value = 0.33319936175468
$ws.Cells.Item(11, 4).Value = 0.2560090365128502
$ws.Cells.Item(11, 5).Value = 0.194010052836056
$ws.Cells.Item(11, 6).Value = 1.1638183204999
$ws.Cells.Item(11, 7).Value = 0.6651295264913699
$ws.Cells.Item(11, 8).Value = 0.6297279495998964
$ws.Cells.Item(11, 9).Value = 0.4280806917633839
$ws.Cells.Item(11, 10).Value = 0.2036191076019378
$ws.Cells.Item(11, 14).Value = 2.017317443662478
$ws.Cells.Item(11, 15).Value = 2.578233826267081

$ws.Cells.Item(12, 2).Value = 0.3408412057247006
$ws.Cells.Item(12, 4).Value = 0.256745331136699
$ws.Cells.Item(12, 5).Value = 0.195178431529726
$ws.Cells.Item(12, 6).Value = 1.174230611703237
$ws.Cells.Item(12, 7).Value = 0.6744288402402958
$ws.Cells.Item(12, 8).Value = 0.6331513654079117
$ws.Cells.Item(12, 9).Value = 0.4316602591565939
$ws.Cells.Item(12, 10).Value = 0.2054348799560159
$ws.Cells.Item(12, 14).Value = 2.048815004509549
$ws.Cells.Item(12, 15).Value = 2.605198180176444

$ws.Cells.Item(13, 2).Value = 0.3391958157679937
$ws.Cells.Item(13, 4).Value = 0.2565860458044114
$ws.Cells.Item(13, 5).Value = 0.1949261616900486
$ws.Cells.Item(13, 6).Value = 1.17198419430666
$ws.Cells.Item(13, 7).Value = 0.6724238672803722
$ws.Cells.Item(13, 8).Value = 0.6324120554357648
$ws.Cells.Item(13, 9).Value = 0.4308881605543888
$ws.Cells.Item(13, 10).Value = 0.2050431063757969
$ws.Cells.Item(13, 14).Value = 2.042030543333738
$ws.Cells.Item(13, 15).Value = 2.599382130016579

$ws.Cells.Item(14, 2).Value = 0.3338282482071406
$ws.Cells.Item(14, 4).Value = 0.2560692914184273
$ws.Cells.Item(14, 5).Value = 0.1941058882925901
$ws.Cells.Item(14, 6).Value = 1.164673168681375
$ws.Cells.Item(14, 7).Value = 0.6658935957561596
$ws.Cells.Item(14, 8).Value = 0.6300086875166357
$ws.Cells.Item(14, 9).Value = 0.4283746547068787
$ws.Cells.Item(14, 10).Value = 0.2037681703561276
$ws.Cells.Item(14, 14).Value = 2.019908357050213
$ws.Cells.Item(14, 15).Value = 2.580448237419887

$ws.Cells.Item(15, 2).Value = 0.3305392456628908
$ws.Cells.Item(15, 4).Value = 0.2557548473068891
$ws.Cells.Item(15, 5).Value = 0.1936053170411043
$ws.Cells.Item(15, 6).Value = 1.160206502899541
$ws.Cells.Item(15, 7).Value = 0.6619000535299335
$ws.Cells.Item(15, 8).Value = 0.6285424593906441
$ws.Cells.Item(15, 9).Value = 0.4268385096919047
$ws.Cells.Item(15, 10).Value = 0.202989326893956
$ws.Cells.Item(15, 14).Value = 2.006360570313291
$ws.Cells.Item(15, 15).Value = 2.568876427918497

$ws.Cells.Item(16, 2).Value = 0.3116732098856971
$ws.Cells.Item(16, 4).Value = 0.253984848216092
$ws.Cells.Item(16, 5).Value = 0.19076540734887
$ws.Cells.Item(16, 6).Value = 1.13478682902678
$ws.Cells.Item(16, 7).Value = 0.6391132448457597
$ws.Cells.Item(16, 8).Value = 0.6202306213541249
$ws.Cells.Item(16, 9).Value = 0.4180888088535895
$ws.Cells.Item(16, 10).Value = 0.1985581258952323
$ws.Cells.Item(16, 14).Value = 1.928767661017872
$ws.Cells.Item(16, 15).Value = 2.502957520220775

$ws.Cells.Item(17, 2).Value = 0.300083899312483
$ws.Cells.Item(17, 4).Value = 0.252927554079335
$ws.Cells.Item(17, 5).Value = 0.1890487917633763
$ws.Cells.Item(17, 6).Value = 1.119350579213645
$ws.Cells.Item(17, 7).Value = 0.6252221158806037
$ws.Cells.Item(17, 8).Value = 0.6152123840539616
$ws.Cells.Item(17, 9).Value = 0.4127692216167418
$ws.Cells.Item(17, 10).Value = 0.1958682884033323
$ws.Cells.Item(17, 14).Value = 1.881209144296406
$ws.Cells.Item(17, 15).Value = 2.462870201680289

$ws.Cells.Item(18, 2).Value = 0.293412510685414
$ws.Cells.Item(18, 4).Value = 0.2523299281001528
$ws.Cells.Item(18, 5).Value = 0.1880708399873683
$ws.Cells.Item(18, 6).Value = 1.110530072250853
$ws.Cells.Item(18, 7).Value = 0.6172645418874652
$ws.Cells.Item(18, 8).Value = 0.6123557350702526
$ws.Cells.Item(18, 9).Value = 0.4097273953839604
$ws.Cells.Item(18, 10).Value = 0.1943316508989312
$ws.Cells.Item(18, 14).Value = 1.853871441567861
$ws.Cells.Item(18, 15).Value = 2.439942247904867

$ws.Cells.Item(19, 2).Value = 0.2911527646049592
$ws.Cells.Item(19, 4).Value = 0.2521293870371721
$ws.Cells.Item(19, 5).Value = 0.1877413368566323
$ws.Cells.Item(19, 6).Value = 1.107553563374182
$ws.Cells.Item(19, 7).Value = 0.614575776360283
$ws.Cells.Item(19, 8).Value = 0.6113936264419806
$ws.Cells.Item(19, 9).Value = 0.4087005711172083
$ws.Cells.Item(19, 10).Value = 0.19381317184191
$ws.Cells.Item(19, 14).Value = 1.844618345498816
$ws.Cells.Item(19, 15).Value = 2.432201422091509

$ws.Cells.Item(20, 2).Value = 0.3013181776856015
$ws.Cells.Item(20, 4).Value = 0.2530390182400026
$ws.Cells.Item(20, 5).Value = 0.1892305555620766
$ws.Cells.Item(20, 6).Value = 1.120987789708067
$ws.Cells.Item(20, 7).Value = 0.6266975129551611
$ws.Cells.Item(20, 8).Value = 0.6157435095531412
$ws.Cells.Item(20, 9).Value = 0.4133336567293782
$ws.Cells.Item(20, 10).Value = 0.1961535407199193
$ws.Cells.Item(20, 14).Value = 1.886270128093628
$ws.Cells.Item(20, 15).Value = 2.467124190622599

$ws.Cells.Item(21, 2).Value = 0.3354050867557419
$ws.Cells.Item(21, 4).Value = 0.256220640706772
$ws.Cells.Item(21, 5).Value = 0.1943464327933526
$ws.Cells.Item(21, 6).Value = 1.166818187367696
$ws.Cells.Item(21, 7).Value = 0.6678103529505393
$ws.Cells.Item(21, 8).Value = 0.6307133845868691
$ws.Cells.Item(21, 9).Value = 0.4291122144186588
$ws.Cells.Item(21, 10).Value = 0.2041422140131317
$ws.Cells.Item(21, 14).Value = 2.026405625114307
$ws.Cells.Item(21, 15).Value = 2.586004214704531

$ws.Cells.Item(22, 2).Value = 0.3576291578720827
$ws.Cells.Item(22, 4).Value = 0.2583932720500997
$ws.Cells.Item(22, 5).Value = 0.1977736478161205
$ws.Cells.Item(22, 6).Value = 1.197287856512247
$ws.Cells.Item(22, 7).Value = 0.6949679831682829
$ws.Cells.Item(22, 8).Value = 0.6407613510679369
$ws.Cells.Item(22, 9).Value = 0.4395793436413342
$ws.Cells.Item(22, 10).Value = 0.2094568644435526
$ws.Cells.Item(22, 14).Value = 2.118115898853603
$ws.Cells.Item(22, 15).Value = 2.664851336729612

$ws.Cells.Item(23, 2).Value = 0.3457728835225566
$ws.Cells.Item(23, 4).Value = 0.2572251772788974
$ws.Cells.Item(23, 5).Value = 0.1959368213278907
$ws.Cells.Item(23, 6).Value = 1.180978317030494
$ws.Cells.Item(23, 7).Value = 0.6804470454142972
$ws.Cells.Item(23, 8).Value = 0.6353743924378534
$ws.Cells.Item(23, 9).Value = 0.4339788619107665
$ws.Cells.Item(23, 10).Value = 0.206611761512363
$ws.Cells.Item(23, 14).Value = 2.069158297356012
$ws.Cells.Item(23, 15).Value = 2.622663646788794

$ws.Cells.Item(24, 2).Value = 0.3007601869816199
$ws.Cells.Item(24, 4).Value = 0.2529885934366405
$ws.Cells.Item(24, 5).Value = 0.1891483522605455
$ws.Cells.Item(24, 6).Value = 1.120247438639538
$ws.Cells.Item(24, 7).Value = 0.6260303967945617
$ws.Cells.Item(24, 8).Value = 0.6155032992347458
$ws.Cells.Item(24, 9).Value = 0.4130784242512391
$ws.Cells.Item(24, 10).Value = 0.1960245476698788
$ws.Cells.Item(24, 14).Value = 1.883982043575145
$ws.Cells.Item(24, 15).Value = 2.465200592235362

$ws.Cells.Item(25, 2).Value = 0.2519980490937712
$ws.Cells.Item(25, 4).Value = 0.2488203805644389
$ws.Cells.Item(25, 5).Value = 0.1821852480318853
$ws.Cells.Item(25, 6).Value = 1.056956987943394
$ws.Cells.Item(25, 7).Value = 0.5685648655344551
$ws.Cells.Item(25, 8).Value = 0.5952049802152146
$ws.Cells.Item(25, 9).Value = 0.391219104399859
$ws.Cells.Item(25, 10).Value = 0.1850052413101935
$ws.Cells.Item(25, 14).Value = 1.684892220676289
$ws.Cells.Item(25, 15).Value = 2.300289579330752
